# Change the household_id field in the "survey" sheet of the household
# member form from an editable, readonly "text" field into a "note" field
# (not editable), per the commit message:
#   "Change the household id in the household members form to be a note
#    (not editable)"
#
# Concretely (row 2 of the "survey" sheet, which defines the household_id
# field):
#   - type (A2):            "text"               -> "note"
#   - name (C2):             "household_id"       -> (cleared)
#   - display.text (D2):     old long description -> "Data for household: {{household_id}}"
#   - readonly (old F2) / comments (old G2) columns are removed entirely
#     (the "readonly" column and its two related long explanatory strings
#     are no longer needed once the field is a plain note).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Update the household_id row to be a "note" field with a templated label.
$ws.Range("A2").Value = "note"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "Data for household: {{household_id}}"

# Drop the now-unused "readonly" (F) and "comments" (G) columns entirely,
# shifting the trailing "hideInContents" column (H) left into column F.
$ws.Range("F1:G1").EntireColumn.Delete()

# Leave the cursor where the author's final edit landed (the
# hideInContents flag for the "note" row, now in column F instead of H).
$ws.Range("F6").Select() | Out-Null
